$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = 0
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 0
$ws.Range("AK2").Value = 0
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = 0
$ws.Range("AN2").Value = 0
$ws.Range("AO2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 0

# Row 3
$ws.Range("AB3").Value = 2.92
$ws.Range("AC3").Value = 8.6
$ws.Range("AD3").Value = 65
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 4.1
$ws.Range("AG3").Value = 12
$ws.Range("AH3").Value = 90
$ws.Range("AJ3").Value = 12.5
$ws.Range("AK3").Value = 46
$ws.Range("AL3").Value = 270
$ws.Range("AN3").Value = 44
$ws.Range("F3").Value = 1.19
$ws.Range("G3").Value = 1.21
$ws.Range("H3").Value = 40
$ws.Range("I3").Value = 55
$ws.Range("J3").Value = 7
$ws.Range("K3").Value = 7.6
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 2.8
$ws.Range("O3").Value = 1.53
$ws.Range("P3").Value = 1.38
$ws.Range("Q3").Value = 3.55
$ws.Range("R3").Value = 1.08
$ws.Range("S3").Value = 11.5
$ws.Range("T3").Value = 3.3
$ws.Range("U3").Value = 1.39
$ws.Range("V3").Value = 1.02
$ws.Range("W3").Value = 5.9
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000

# Row 4
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 6.6
$ws.Range("AD4").Value = 15
$ws.Range("AE4").Value = 65
$ws.Range("AF4").Value = 7.6
$ws.Range("AG4").Value = 7.2
$ws.Range("AH4").Value = 14.5
$ws.Range("AI4").Value = 60
$ws.Range("AJ4").Value = 17
$ws.Range("AK4").Value = 18
$ws.Range("AL4").Value = 36
$ws.Range("AM4").Value = 110
$ws.Range("AN4").Value = 20
$ws.Range("AO4").Value = 110
$ws.Range("F4").Value = 1.57
$ws.Range("G4").Value = 1.58
$ws.Range("H4").Value = 8
$ws.Range("I4").Value = 8.4
$ws.Range("J4").Value = 4.1
$ws.Range("K4").Value = 4.3
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 6.4
$ws.Range("O4").Value = 1.17
$ws.Range("P4").Value = 2.28
$ws.Range("Q4").Value = 1.77
$ws.Range("R4").Value = 1.39
$ws.Range("S4").Value = 3.5
$ws.Range("T4").Value = 1.56
$ws.Range("U4").Value = 2.62
$ws.Range("V4").Value = 1.13
$ws.Range("W4").Value = 2.72
$ws.Range("X4").Value = 1000
$ws.Range("Y4").Value = 1000
$ws.Range("Z4").Value = 1000

# Row 5
$ws.Range("AA5").Value = 170
$ws.Range("AB5").Value = 10.5
$ws.Range("AC5").Value = 9.6
$ws.Range("AD5").Value = 23
$ws.Range("AE5").Value = 70
$ws.Range("AF5").Value = 11.5
$ws.Range("AG5").Value = 9.4
$ws.Range("AH5").Value = 18
$ws.Range("AI5").Value = 70
$ws.Range("AK5").Value = 16
$ws.Range("AL5").Value = 32
$ws.Range("AM5").Value = 110
$ws.Range("AN5").Value = 9
$ws.Range("AO5").Value = 80
$ws.Range("F5").Value = 1.71
$ws.Range("G5").Value = 1.72
$ws.Range("H5").Value = 5.7
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 4.2
$ws.Range("L5").Value = 1.47
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 5.3
$ws.Range("O5").Value = 1.22
$ws.Range("P5").Value = 2.42
$ws.Range("Q5").Value = 1.66
$ws.Range("R5").Value = 1.58
$ws.Range("S5").Value = 2.66
$ws.Range("T5").Value = 1.66
$ws.Range("U5").Value = 2.34
$ws.Range("V5").Value = 1.19
$ws.Range("W5").Value = 2.44
$ws.Range("X5").Value = 23
$ws.Range("Y5").Value = 30
$ws.Range("Z5").Value = 55

# Row 6
$ws.Range("AA6").Value = 25
$ws.Range("AB6").Value = 46
$ws.Range("AD6").Value = 11.5
$ws.Range("AE6").Value = 26
$ws.Range("AF6").Value = 48
$ws.Range("AG6").Value = 27
$ws.Range("AH6").Value = 32
$ws.Range("AI6").Value = 75
$ws.Range("AJ6").Value = 230
$ws.Range("AK6").Value = 160
$ws.Range("AL6").Value = 190
$ws.Range("AM6").Value = 340
$ws.Range("AN6").Value = 340
$ws.Range("AO6").Value = 21
$ws.Range("C6").Value = "14:30:02"
$ws.Range("G6").Value = 5.9
$ws.Range("H6").Value = 1.83
$ws.Range("I6").Value = 1.87
$ws.Range("J6").Value = 3.5
$ws.Range("K6").Value = 3.6
$ws.Range("L6").Value = 1.83
$ws.Range("M6").Value = 1.11
$ws.Range("N6").Value = 2.7
$ws.Range("O6").Value = 1.55
$ws.Range("P6").Value = 1.57
$ws.Range("Q6").Value = 2.72
$ws.Range("R6").Value = 1.2
$ws.Range("S6").Value = 5.6
$ws.Range("T6").Value = 2.42
$ws.Range("U6").Value = 1.64
$ws.Range("V6").Value = 2.2
$ws.Range("W6").Value = 1.19
$ws.Range("X6").Value = 9.6
$ws.Range("Y6").Value = 6.2
$ws.Range("Z6").Value = 8.800000000000001

# Row 7
$ws.Range("AC7").Value = 8.4
$ws.Range("AE7").Value = 1000
$ws.Range("AF7").Value = 9.4
$ws.Range("AJ7").Value = 20
$ws.Range("AK7").Value = 26
$ws.Range("AL7").Value = 95
$ws.Range("AN7").Value = 22
$ws.Range("F7").Value = 1.74
$ws.Range("G7").Value = 1.82
$ws.Range("H7").Value = 6.2
$ws.Range("I7").Value = 7.2
$ws.Range("J7").Value = 3.35
$ws.Range("K7").Value = 3.6
$ws.Range("L7").Value = 1.58
$ws.Range("N7").Value = 2.68
$ws.Range("O7").Value = 1.56
$ws.Range("Q7").Value = 2.7
$ws.Range("S7").Value = 5.6
$ws.Range("T7").Value = 2.34
$ws.Range("V7").Value = 1.17
$ws.Range("W7").Value = 2.2
$ws.Range("X7").Value = 9
$ws.Range("Y7").Value = 16
$ws.Range("Z7").Value = 55

# Row 8
$ws.Range("AA8").Value = 400
$ws.Range("AB8").Value = 6
$ws.Range("AC8").Value = 9.6
$ws.Range("AD8").Value = 36
$ws.Range("AE8").Value = 190
$ws.Range("AI8").Value = 180
$ws.Range("AK8").Value = 19
$ws.Range("AN8").Value = 11
$ws.Range("AO8").Value = 330
$ws.Range("F8").Value = 1.51
$ws.Range("G8").Value = 1.52
$ws.Range("H8").Value = 9
$ws.Range("I8").Value = 9.199999999999999
$ws.Range("J8").Value = 4.3
$ws.Range("K8").Value = 4.4
$ws.Range("N8").Value = 3.25
$ws.Range("P8").Value = 1.75
$ws.Range("Q8").Value = 2.3
$ws.Range("S8").Value = 4.5
$ws.Range("V8").Value = 1.12
$ws.Range("W8").Value = 2.92
$ws.Range("Y8").Value = 22
$ws.Range("Z8").Value = 75

# Row 9
$ws.Range("AC9").Value = 9.800000000000001
$ws.Range("AH9").Value = 21
$ws.Range("AJ9").Value = 15.5
$ws.Range("AK9").Value = 15
$ws.Range("AL9").Value = 30
$ws.Range("AM9").Value = 100
$ws.Range("AN9").Value = 8.6
$ws.Range("AO9").Value = 90
$ws.Range("F9").Value = 1.61
$ws.Range("G9").Value = 1.62
$ws.Range("J9").Value = 4.5
$ws.Range("K9").Value = 4.6
$ws.Range("L9").Value = 1.4
$ws.Range("N9").Value = 4.3
$ws.Range("O9").Value = 1.29
$ws.Range("P9").Value = 2.14
$ws.Range("Q9").Value = 1.86
$ws.Range("R9").Value = 1.44
$ws.Range("S9").Value = 3.15
$ws.Range("T9").Value = 1.92
$ws.Range("W9").Value = 2.6

# Row 10
$ws.Range("AA10").Value = 900
$ws.Range("AB10").Value = 8.4
$ws.Range("AE10").Value = 60
$ws.Range("AF10").Value = 13.5
$ws.Range("AG10").Value = 11
$ws.Range("AH10").Value = 21
$ws.Range("AI10").Value = 130
$ws.Range("AK10").Value = 26
$ws.Range("AL10").Value = 48
$ws.Range("AN10").Value = 23
$ws.Range("F10").Value = 2.2
$ws.Range("H10").Value = 3.85
$ws.Range("I10").Value = 4.2
$ws.Range("J10").Value = 3.2
$ws.Range("K10").Value = 3.35
$ws.Range("M10").Value = 1.1
$ws.Range("O10").Value = 1.44
$ws.Range("P10").Value = 1.7
$ws.Range("S10").Value = 4.6
$ws.Range("T10").Value = 1.95
$ws.Range("U10").Value = 2
$ws.Range("V10").Value = 1.31
$ws.Range("W10").Value = 1.8
$ws.Range("X10").Value = 10.5
$ws.Range("Y10").Value = 16
$ws.Range("Z10").Value = 29
